# Users schema sheet: fix ColumnID casing for the LDAP/agent/company fields
# (companyEmail -> companyemail, mobileEmail -> mobileemail, agentCD -> agentcd,
#  agentStrDay -> agentstrday, agentEndDay -> agentendday, companyTel -> companytel,
#  deleteReason -> deletereason, DN -> dn) to match the lower-cased column names
# used elsewhere in the schema definition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "dn"
$ws.Range("A12").Value = "companyemail"
$ws.Range("A13").Value = "mobileemail"
$ws.Range("A14").Value = "agentcd"
$ws.Range("A15").Value = "agentstrday"
$ws.Range("A16").Value = "agentendday"
$ws.Range("A17").Value = "companytel"
$ws.Range("A21").Value = "deletereason"

# The workbook was re-saved from a Japanese-locale Excel, which localizes the
# built-in cell style names and the theme name as a side effect of opening &
# saving the file. Re-apply that localization here.
$styleMap = @{
  "Normal" = "標準"
  "Title" = "タイトル"
  "Heading 1" = "見出し 1"
  "Heading 2" = "見出し 2"
  "Heading 3" = "見出し 3"
  "Heading 4" = "見出し 4"
  "Good" = "良い"
  "Bad" = "悪い"
  "Neutral" = "どちらでもない"
  "Input" = "入力"
  "Output" = "出力"
  "Calculation" = "計算"
  "Linked Cell" = "リンク セル"
  "Check Cell" = "チェック セル"
  "Warning Text" = "警告文"
  "Note" = "メモ"
  "Explanatory Text" = "説明文"
  "Total" = "集計"
  "Accent1" = "アクセント 1"
  "20% - Accent1" = "20% - アクセント 1"
  "40% - Accent1" = "40% - アクセント 1"
  "60% - Accent1" = "60% - アクセント 1"
  "Accent2" = "アクセント 2"
  "20% - Accent2" = "20% - アクセント 2"
  "40% - Accent2" = "40% - アクセント 2"
  "60% - Accent2" = "60% - アクセント 2"
  "Accent3" = "アクセント 3"
  "20% - Accent3" = "20% - アクセント 3"
  "40% - Accent3" = "40% - アクセント 3"
  "60% - Accent3" = "60% - アクセント 3"
  "Accent4" = "アクセント 4"
  "20% - Accent4" = "20% - アクセント 4"
  "40% - Accent4" = "40% - アクセント 4"
  "60% - Accent4" = "60% - アクセント 4"
  "Accent5" = "アクセント 5"
  "20% - Accent5" = "20% - アクセント 5"
  "40% - Accent5" = "40% - アクセント 5"
  "60% - Accent5" = "60% - アクセント 5"
  "Accent6" = "アクセント 6"
  "20% - Accent6" = "20% - アクセント 6"
  "40% - Accent6" = "40% - アクセント 6"
  "60% - Accent6" = "60% - アクセント 6"
}

for ($i = 1; $i -le $wb.Styles.Count; $i++) {
  $s = $wb.Styles.Item($i)
  $oldName = $s.Name
  if ($styleMap.ContainsKey($oldName)) {
    $s.Name = $styleMap[$oldName]
  }
}

try {
  $wb.Theme.Name = "Office ​​テーマ"
} catch {
}

# Cursor was left on A27 (just below the table) when the file was saved.
$ws.Range("A27").Select()
